$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the xpath value in B2 to reflect the SOAP response path instead of the request path
$ws.Range("B2").Value = "/Envelope/Body/getActorsByIdResponse/actor/actor_id"

# Add the new expected value cell C2 (second SOAP test)
$ws.Range("C2").Value = 5

# Resize column B to fit the longer xpath text
$ws.Columns("B").AutoFit()

# Update the selection to the full second row, matching the state after the edits
[void]$ws.Range("A2:XFD2").Select()
